$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "losstan" column header
$ws.Range("I1").Value = "losstan"

# Fill in the loss tangent values for each cable type row
$ws.Range("I2").Value = 0.00007
$ws.Range("I3").Value = 0.0003
$ws.Range("I4").Value = 0.0002
$ws.Range("I5").Value = 0.00007

# Apply the same number format (scientific notation) used by column H to the new column
$ws.Range("I2:I5").NumberFormat = "0.00E+00"

# Update the active selection/view to reflect the new cell of interest
$ws.Range("I5").Select()
